$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Actual Duration (F) and Percent Complete (G) for rows 11-14
# Row 11: API Endpoints
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 1

# Row 12: Database Completion
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 1

# Row 13: Frontend Development
$ws.Range("F13").Value = 7
$ws.Range("G13").Value = 1

# Row 14: Frontend Polish
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 0.66
